$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 46052
$ws.Range("D8").Value = 157.95
$ws.Range("E8").Value = 149.37
$ws.Range("F8").Value = 159.37
$ws.Range("G8").Value = 149.39
$ws.Range("A9").Value = 46052
$ws.Range("D9").Value = 157.95
$ws.Range("E9").Value = 149.37
$ws.Range("F9").Value = 159.37
$ws.Range("G9").Value = 149.39
$ws.Range("A10").Value = 46052
$ws.Range("D10").Value = 158.77
$ws.Range("E10").Value = 151.17
$ws.Range("F10").Value = 161.17
$ws.Range("G10").Value = 151.57
$ws.Range("A11").Value = 46051
$ws.Range("D11").Value = 158.25
$ws.Range("E11").Value = 149.89
$ws.Range("F11").Value = 159.89
$ws.Range("G11").Value = 149.91
$ws.Range("A12").Value = 46051
$ws.Range("D12").Value = 158.25
$ws.Range("E12").Value = 149.89
$ws.Range("F12").Value = 159.89
$ws.Range("G12").Value = 149.91
$ws.Range("A13").Value = 46051
$ws.Range("D13").Value = 159.12
$ws.Range("E13").Value = 151.73
$ws.Range("F13").Value = 161.73
$ws.Range("G13").Value = 152.13
$ws.Range("A17").Value = 46052
$ws.Range("D17").Value = 162.3
$ws.Range("E17").Value = 153.79
$ws.Range("F17").Value = 163.79
$ws.Range("A18").Value = 46051
$ws.Range("D18").Value = 162.61
$ws.Range("E18").Value = 154.32
$ws.Range("F18").Value = 164.32
$ws.Range("A22").Value = 46052
$ws.Range("D22").Value = 158.91
$ws.Range("E22").Value = 150.94
$ws.Range("F22").Value = 160.54
$ws.Range("G22").Value = 152.02
$ws.Range("A23").Value = 46052
$ws.Range("D23").Value = 163.54
$ws.Range("E23").Value = 156.51
$ws.Range("F23").Value = 166.51
$ws.Range("A24").Value = 46052
$ws.Range("D24").Value = 163.69
$ws.Range("E24").Value = 157.18
$ws.Range("F24").Value = 167.18
$ws.Range("A25").Value = 46052
$ws.Range("D25").Value = 163.69
$ws.Range("E25").Value = 156.71
$ws.Range("F25").Value = 166.71
$ws.Range("G25").Value = 156.84
$ws.Range("A26").Value = 46052
$ws.Range("D26").Value = 163.28
$ws.Range("E26").Value = 158.31
$ws.Range("F26").Value = 168.31
$ws.Range("A27").Value = 46051
$ws.Range("D27").Value = 159.21
$ws.Range("E27").Value = 151.46
$ws.Range("F27").Value = 161.06
$ws.Range("G27").Value = 152.54
$ws.Range("A28").Value = 46051
$ws.Range("D28").Value = 163.89
$ws.Range("E28").Value = 157.07
$ws.Range("F28").Value = 167.07
$ws.Range("A29").Value = 46051
$ws.Range("D29").Value = 164.04
$ws.Range("E29").Value = 157.75
$ws.Range("F29").Value = 167.75
$ws.Range("A30").Value = 46051
$ws.Range("D30").Value = 164.03
$ws.Range("E30").Value = 157.29
$ws.Range("F30").Value = 167.29
$ws.Range("G30").Value = 157.41
$ws.Range("A31").Value = 46051
$ws.Range("D31").Value = 163.62
$ws.Range("E31").Value = 158.89
$ws.Range("F31").Value = 168.89
$ws.Range("A35").Value = 46052
$ws.Range("D35").Value = 157.69
$ws.Range("E35").Value = 148.19
$ws.Range("F35").Value = 157.19
$ws.Range("A36").Value = 46051
$ws.Range("D36").Value = 158.02
$ws.Range("E36").Value = 148.74
$ws.Range("F36").Value = 157.74
$ws.Range("A40").Value = 46052
$ws.Range("D40").Value = 163.68
$ws.Range("E40").Value = 156.15
$ws.Range("F40").Value = 166.15
$ws.Range("A41").Value = 46052
$ws.Range("D41").Value = 163.4
$ws.Range("E41").Value = 156.57
$ws.Range("F41").Value = 166.57
$ws.Range("A42").Value = 46051
$ws.Range("D42").Value = 163.49
$ws.Range("E42").Value = 156.32
$ws.Range("F42").Value = 166.32
$ws.Range("A43").Value = 46051
$ws.Range("D43").Value = 163.21
$ws.Range("E43").Value = 156.74
$ws.Range("F43").Value = 166.74
$ws.Range("A47").Value = 46052
$ws.Range("D47").Value = 157.57
$ws.Range("E47").Value = 150.14
$ws.Range("F47").Value = 160.14
$ws.Range("A48").Value = 46052
$ws.Range("D48").Value = 157.17
$ws.Range("E48").Value = 150.06
$ws.Range("F48").Value = 160.06
$ws.Range("A49").Value = 46051
$ws.Range("D49").Value = 157.86
$ws.Range("E49").Value = 150.44
$ws.Range("F49").Value = 160.44
$ws.Range("A50").Value = 46051
$ws.Range("D50").Value = 157.45
$ws.Range("E50").Value = 150.36
$ws.Range("F50").Value = 160.36
$ws.Range("A54").Value = 46052
$ws.Range("D54").Value = 172.13
$ws.Range("E54").Value = 164.15
$ws.Range("F54").Value = 174.15
$ws.Range("A55").Value = 46052
$ws.Range("D55").Value = 164.78
$ws.Range("E55").Value = 162.48
$ws.Range("F55").Value = 172.48
$ws.Range("A56").Value = 46052
$ws.Range("D56").Value = 161.72
$ws.Range("A57").Value = 46052
$ws.Range("D57").Value = 162.21
$ws.Range("E57").Value = 156.9
$ws.Range("A58").Value = 46052
$ws.Range("D58").Value = 157.98
$ws.Range("E58").Value = 152.8
$ws.Range("F58").Value = 162.8
$ws.Range("A59").Value = 46052
$ws.Range("D59").Value = 164.65
$ws.Range("E59").Value = 162.27
$ws.Range("A60").Value = 46051
$ws.Range("D60").Value = 172.51
$ws.Range("E60").Value = 164.78
$ws.Range("F60").Value = 174.78
$ws.Range("A61").Value = 46051
$ws.Range("D61").Value = 165.19
$ws.Range("E61").Value = 163.01
$ws.Range("F61").Value = 173.01
$ws.Range("A62").Value = 46051
$ws.Range("D62").Value = 162.06
$ws.Range("A63").Value = 46051
$ws.Range("D63").Value = 162.52
$ws.Range("E63").Value = 157.44
$ws.Range("A64").Value = 46051
$ws.Range("D64").Value = 158.29
$ws.Range("E64").Value = 153.33
$ws.Range("F64").Value = 163.33
$ws.Range("A65").Value = 46051
$ws.Range("D65").Value = 164.95
$ws.Range("E65").Value = 162.89
